$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja3")
$ws.Range("A32").Value = 28
$ws.Range("B32").Value = "obtener_nombre_usuario()"
$ws.Range("C32").Value = "genera el nombre de usuario que se va a registrar en los campos CREATED_BY de cada tabla afectada"
$ws.Range("D32").Value = "OK"
